$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: shift D1/E1 left into C1/D1, and move old C1 value (max) to E1
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Row 2: C2 becomes the taxonomy string (same as D2), D2 unchanged, E2 becomes numeric 1
$ws.Range("C2").Value = "o__Elusimicrobiales"
$ws.Range("E2").Value = 1

# Row 3: same pattern
$ws.Range("C3").Value = "o__Elusimicrobiales"
$ws.Range("E3").Value = 1
